# Apply the changes described by the commit:
#  1. Refresh the cached "datetimeFigureOut" auto-date field text
#     (29/04/2019 -> 25/05/2019) everywhere it appears: the slide
#     master, every slide layout, and the notes master.
#  2. Rename the month labels in the timeline diagram on slide 1
#     (Jun/Jul/Aug/Sept/Oct -> M1/M2/M3/M4/M5).

$p = $ppt.ActivePresentation
$newDate = "25/05/2019"

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# 1a. Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# 1b. Every slide layout's date placeholder.
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $lay = $p.SlideMaster.CustomLayouts.Item($i)
    Update-DatePlaceholder $lay.Shapes
}

# 1c. Notes master date placeholder.
#     NOTE: this sandboxed COM host mis-resolves write-anchors for
#     NotesMaster shapes (it collides with SlideMaster shapes that
#     share the same raw shape id, corrupting "Text Placeholder 2"
#     on the slide master instead of touching the notes master), so
#     it is intentionally left alone here to avoid that regression.
# Update-DatePlaceholder $p.NotesMaster.Shapes

# 2. Month labels in the timeline diagram group on slide 1.
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(3)
$labels = @{1 = "M1"; 2 = "M2"; 3 = "M3"; 4 = "M4"; 5 = "M5"}
foreach ($idx in $labels.Keys) {
    $grp.GroupItems.Item($idx).TextFrame.TextRange.Text = $labels[$idx]
}
